# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Golem_Profits Leve-profit workbook.
# Source: unified OOXML diff (chore: update Sheets via scheduled runner)

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 707.7857
$ws.Range("I19").Value = 397
$ws.Range("K19").Value = 397
$ws.Range("M19").Value = -222
$ws.Range("H28").Value = 1901.6364
$ws.Range("I28").Value = 1839.2222
$ws.Range("J28").Value = 2182.5
$ws.Range("K28").Value = 1839.2222
$ws.Range("L28").Value = 2182.5
$ws.Range("M28").Value = -1354.2222
$ws.Range("N28").Value = -3152.5
$ws.Range("H58").Value = 1167.75
$ws.Range("I58").Value = 481.33334
$ws.Range("J58").Value = 1579.6
$ws.Range("K58").Value = 1444.00002
$ws.Range("L58").Value = 4738.799999999999
$ws.Range("M58").Value = -1294.00002
$ws.Range("N58").Value = -5038.799999999999
$ws.Range("H96").Value = 633.3333
$ws.Range("I96").Value = 118.666664
$ws.Range("K96").Value = 355.999992
$ws.Range("M96").Value = 1017.000008
$ws.Range("H113").Value = 18341.715
$ws.Range("I113").Value = 19732.166
$ws.Range("K113").Value = 19732.166
$ws.Range("M113").Value = -16478.166
$ws.Range("H135").Value = 4507.75
$ws.Range("I135").Value = 4507.75
$ws.Range("K135").Value = 40569.75
$ws.Range("M135").Value = -38034.75

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1362.25
$ws.Range("I122").Value = 1233
$ws.Range("K122").Value = 3699
$ws.Range("M122").Value = -1249
$ws.Range("H132").Value = 2587.4285
$ws.Range("I132").Value = 1028
$ws.Range("K132").Value = 3084
$ws.Range("M132").Value = -554
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2545.4546
$ws.Range("I94").Value = 1980
$ws.Range("J94").Value = 3016.6667
$ws.Range("K94").Value = 1980
$ws.Range("L94").Value = 3016.6667
$ws.Range("M94").Value = -1529
$ws.Range("N94").Value = -3918.6667

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 839.7143
$ws.Range("J22").Value = 1147
$ws.Range("L22").Value = 1147
$ws.Range("N22").Value = -1847

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 434.33334
$ws.Range("I18").Value = 432.2
$ws.Range("K18").Value = 1296.6
$ws.Range("M18").Value = -1127.6
$ws.Range("H63").Value = 2945.7778
$ws.Range("I63").Value = 1768.6666
$ws.Range("K63").Value = 5305.9998
$ws.Range("M63").Value = -4556.9998
$ws.Range("H66").Value = 2945.7778
$ws.Range("I66").Value = 1768.6666
$ws.Range("K66").Value = 15917.9994
$ws.Range("M66").Value = -12173.9994
$ws.Range("H70").Value = 9749.25
$ws.Range("I70").Value = 9749.25
$ws.Range("K70").Value = 29247.75
$ws.Range("M70").Value = -28932.75
$ws.Range("H73").Value = 9749.25
$ws.Range("I73").Value = 9749.25
$ws.Range("K73").Value = 29247.75
$ws.Range("M73").Value = -28155.75
$ws.Range("H87").Value = 250
$ws.Range("I87").Value = 250
$ws.Range("K87").Value = 750
$ws.Range("M87").Value = 498
$ws.Range("H90").Value = 250
$ws.Range("I90").Value = 250
$ws.Range("K90").Value = 2250
$ws.Range("M90").Value = 3990
$ws.Range("H132").Value = 9560
$ws.Range("I132").Value = 7333.3335
$ws.Range("K132").Value = 66000.0015
$ws.Range("M132").Value = -63470.0015
$ws.Range("H140").Value = 369.75
$ws.Range("I140").Value = 369.75
$ws.Range("K140").Value = 1109.25
$ws.Range("M140").Value = 4070.75

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H122").Value = 14999
$ws.Range("I122").Value = 14999
$ws.Range("K122").Value = 44997
$ws.Range("M122").Value = -42547

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 933
$ws.Range("I7").Value = 933
$ws.Range("K7").Value = 933
$ws.Range("M7").Value = -821
$ws.Range("H61").Value = 1915.1111
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3404
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 1728.4
$ws.Range("I82").Value = 1771
$ws.Range("J82").Value = 1700
$ws.Range("K82").Value = 1771
$ws.Range("L82").Value = 1700
$ws.Range("M82").Value = -1410
$ws.Range("N82").Value = -2422
$ws.Range("H85").Value = 1728.4
$ws.Range("I85").Value = 1771
$ws.Range("J85").Value = 1700
$ws.Range("K85").Value = 1771
$ws.Range("L85").Value = 1700
$ws.Range("M85").Value = -523
$ws.Range("N85").Value = -4196
$ws.Range("H113").Value = 1915.1111
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H126").Value = 933
$ws.Range("I126").Value = 933
$ws.Range("K126").Value = 2799
$ws.Range("M126").Value = -329

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H113").Value = 14197.5
$ws.Range("I113").Value = 390
$ws.Range("J113").Value = 18800
$ws.Range("K113").Value = 1170
$ws.Range("L113").Value = 56400
$ws.Range("M113").Value = 1000
$ws.Range("N113").Value = -60740
$ws.Range("H132").Value = 3941
$ws.Range("I132").Value = 3950
$ws.Range("J132").Value = 3935
$ws.Range("K132").Value = 11850
$ws.Range("L132").Value = 11805
$ws.Range("M132").Value = -9320
